$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Majorelle Magdy, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad"
$ws.Range("G3").Value = "Dr. Eman Tantawi, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Asmaa Reda"
$ws.Range("G4").Value = "Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Asmaa Reda"
$ws.Range("G5").Value = "Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Nesma, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Mohammad El-Tanany"
$ws.Range("G6").Value = "Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Nahla Nagiub"
$ws.Range("G7").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Asmaa Reda"
$ws.Range("G8").Value = "Dr. Majorelle Magdy, Dr. Eman Tantawi, Dr. Manar Montaser, Administrator, Dr. Asmaa Reda, Dr. Shimaa Ahmad Mekki"
$ws.Range("G9").Value = "Dr. Majorelle Magdy, Dr. Menna tuâ€™Allah Medhat, Dr. Manar Montaser, Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Asmaa Reda"
$ws.Range("G10").Value = "Dr. Gehan Adel, Dr. Alshimaa Atef, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Sara Wael, Dr. Heba Mahmoud Ali, Dr. Shimaa Ahmad Mekki"
$ws.Range("G11").Value = "Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Asmaa Reda"
$ws.Range("G13").Value = "Dr. Mariam Nour El-Din, Dr. Shimaa Ashraf, D Wessam Atef, Dr. Omnia Mohammad, Dr. Safa Hany"
$ws.Range("G14").Value = "Dr. Shimaa Ashraf, Dr. Safa Hany"
$ws.Range("G17").Value = "Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Eman M. Abo-Sakaya, Dr. Sarah Abdelmohsen, Dr. Arwa Al-Sayed, Dr. Dina Adel, Dr. Marwa Mustafa"
$ws.Range("G24").Value = "Dr. Monica, Dr. Marina Atef, Dr. Wafaa Ebida, Dr. Youstina Magdy, Dr. Yasmin, Dr. Maryam Ashraf, Dr. Salma Hassan, Dr. Aya Emad, Dr. Remon, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah"
$ws.Range("G25").Value = "Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Marina Atef, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Remon, Dr. Ola Abd Al-Fattah"
$ws.Range("G27").Value = "Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Yasmin, Dr. Salma Hassan, Dr. Eman Mohammad Al, Dr. Remon, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah"
$ws.Range("G28").Value = "Dr. Eman Samir Gabry, Dr. Nardine, Dr. Abdullah El-Agrody, Dr. Wafaa Ebida, Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Remon, Dr. Neveen Nashaat"
$ws.Range("G29").Value = "Dr. Eman Samir Gabry, Dr. Monica, Dr. Naema Gomaa, Dr. Remon, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah"
$ws.Range("G30").Value = "Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Asmaa Reda"
$ws.Range("G31").Value = "Dr. Eman Tantawi, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Asmaa Reda"
$ws.Range("G32").Value = "Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Veronia Rafat, Dr. Asmaa Reda"
$ws.Range("G33").Value = "Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Nesma, Dr. Hanan Ragab, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Mohammad El-Tanany"
$ws.Range("G34").Value = "Dr. Nourhan Mahmoud, Dr. Eman Tantawi, Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Asmaa Reda, Dr. Nahla Nagiub"
$ws.Range("G35").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Asmaa Reda"
$ws.Range("G36").Value = "Dr. Majorelle Magdy, Dr. Eman Tantawi, Dr. Manar Montaser, Administrator, Dr. Asmaa Reda, Dr. Shimaa Ahmad Mekki"
$ws.Range("G37").Value = "Dr. Majorelle Magdy, Dr. Menna tuâ€™Allah Medhat, Dr. Manar Montaser, Dr. Gehan Adel, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Rana Abo-Zaid, Dr. Asmaa Reda"
$ws.Range("G38").Value = "Dr. Gehan Adel, Dr. Alshimaa Atef, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad, Dr. Sara Wael, Dr. Heba Mahmoud Ali, Dr. Shimaa Ahmad Mekki"
$ws.Range("G39").Value = "Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Asmaa Reda"
$ws.Range("G41").Value = "Dr. Mariam Nour El-Din, Dr. Shimaa Ashraf, D Wessam Atef, Dr. Omnia Mohammad, Dr. Safa Hany"
$ws.Range("G42").Value = "Dr. Shimaa Ashraf, Dr. Safa Hany"
$ws.Range("G45").Value = "Dr. Basma Hamed, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Eman M. Abo-Sakaya, Dr. Sarah Abdelmohsen, Dr. Arwa Al-Sayed, Dr. Dina Adel, Dr. Marwa Mustafa"
$ws.Range("G52").Value = "Dr. Monica, Dr. Marina Atef, Dr. Wafaa Ebida, Dr. Youstina Magdy, Dr. Yasmin, Dr. Maryam Ashraf, Dr. Salma Hassan, Dr. Aya Emad, Dr. Remon, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah"
$ws.Range("G53").Value = "Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Marina Atef, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Remon, Dr. Ola Abd Al-Fattah"
$ws.Range("G55").Value = "Dr. Eman Samir Gabry, Dr. Wafaa Ebida, Dr. Yasmin, Dr. Salma Hassan, Dr. Eman Mohammad Al, Dr. Remon, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah"
$ws.Range("G56").Value = "Dr. Eman Samir Gabry, Dr. Nardine, Dr. Abdullah El-Agrody, Dr. Wafaa Ebida, Dr. Aya Hanafy, Dr. Salma Hassan, Dr. Remon, Dr. Neveen Nashaat"
$ws.Range("G57").Value = "Dr. Eman Samir Gabry, Dr. Monica, Dr. Naema Gomaa, Dr. Remon, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah"
